# "Fixed issues with Edit Order"
#
# The state column (F) held shouty, all-caps state names ("FLORIDA",
# "MICHIGAN") together with a handful of bogus "C@bi$ush5" mailto
# hyperlinks that had been accidentally stamped onto F2:F5. This script:
#   1) normalises the state values to proper case (Florida / Michigan),
#   2) drops the bogus F2:F5 hyperlinks while keeping the legitimate
#      C2:C6 email hyperlinks (re-added so rId numbering stays compact),
#   3) restores the original cell formatting on the email/state cells so
#      only the hyperlink *relationships* changed, not their look, and
#   4) leaves the selection where the user's editing session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. stash the current formatting of a "bordered hyperlink" cell (C2)
#        and a "plain hyperlink" cell (C6) in scratch cells so we can
#        restore it after Hyperlinks.Add() re-stamps a new style on them.
$ws.Range("C2").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("C6").Copy()
$ws.Range("J2").PasteSpecial(-4122)

# --- 2. normalise the state column values (ALL CAPS -> Proper case)
$ws.Range("F2").Value = "Florida"
$ws.Range("F3").Value = "Florida"
$ws.Range("F4").Value = "Florida"
$ws.Range("F5").Value = "Florida"
$ws.Range("F6").Value = "Michigan"

# --- 3. rebuild the hyperlinks collection, keeping only the genuine
#        C-column email links (this removes the stray F2:F5 links).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:neohost1a@test.com", "", "")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:neocohost1a@test.com", "", "")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:neoguest1ab@test.com", "", "")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:neoguest2ab@test.com", "", "")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:michigan@na.com", "", "")

# --- 4. restore the original look of the email cells (Hyperlinks.Add
#        stamps its own style; put the pre-existing formatting back).
$ws.Range("J1").Copy()
$ws.Range("C2:C5").PasteSpecial(-4122)
$ws.Range("J2").Copy()
$ws.Range("C6").PasteSpecial(-4122)

# the now-unlinked F2:F5 cells should go back to plain/unstyled cells
$ws.Range("F2:F5").Style = "Normal"

# tidy up the scratch cells used to stash formatting
$ws.Range("J1:J2").Clear()

# --- 5. restore the end-user's selection
$null = $ws.Range("E15").Select()
